$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "95.919.67"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +4.16%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.087.15"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.31%  "

# Row 4
$ws.Range("E4").Value = "  -0.01%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "236.40"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.10%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "602.84"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.33%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.10"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.00%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.377"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.64%  "

# Row 9
$ws.Range("E9").Value = "  +0.01%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "3.082.66"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.38%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.777"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.02%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.195"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.53%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "95.221.76"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +3.62%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000235"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.93%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "33.18"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.12%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.27"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.11%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.659.87"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.54%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.074.59"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.48%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.48"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -8.38%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.16"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.77%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "452.06"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.95%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.56"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.20%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.0000188"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -5.27%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "8.59"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -4.91%  "

# Row 25
$ws.Range("E25").Value = "  -2.12%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "84.65"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.39%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.45"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.61%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "3.251.92"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.22%  "

# Row 29
$ws.Range("E29").Value = "  -0.04%  "

# Row 30
$ws.Range("B30").Value = "Cronos"
$ws.Range("C30").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.178"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.82%  "

# Row 31
$ws.Range("B31").Value = "Stellar"
$ws.Range("C31").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.238"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +3.74%  "

# Row 32
$ws.Range("B32").Value = "Hedera"
$ws.Range("C32").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.127"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.36%  "

# Row 33
$ws.Range("B33").Value = "InternetComputer(DFINITY)"
$ws.Range("C33").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "8.80"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.14%  "

# Row 34
$ws.Range("B34").Value = "EthereumClassic"
$ws.Range("C34").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "25.41"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.02%  "

# Row 35
$ws.Range("B35").Value = "Binance-PegBSC-USD"
$ws.Range("C35").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.819"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -18.02%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "7.23"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -8.31%  "

# Row 37
$ws.Range("E37").Value = "  -4.72%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "24.13"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.21%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "477.86"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.27%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.82"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.99%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.427"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.12%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.64"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -6.93%  "

# Row 43
$ws.Range("E43").Value = "  -4.52%  "

# Row 44
$ws.Range("E44").Value = "  +0.04%  "

# Row 45
$ws.Range("E45").Value = "  -5.90%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "161.56"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.12%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.669"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.26%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.83"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.74%  "

# Row 49
$ws.Range("E49").Value = "  +13.13%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "43.70"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.16%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.998"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.06%  "
